$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Formula = "=A2/B2"
$ws.Range("D13").Formula = "=B13/B2"
$ws.Range("E13").Formula = "=D13"
$ws.Range("D17").Select()
